$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: swap the F:V (cols 6-22) content of two rows in the odds table.
# Columns A-E (index/pais/torneio/temporada/data_partida) stay tied to the
# row position; F..V (home/away teams, odds, timestamps, url) are what the
# source script re-ordered.
function Swap-RowOdds($r1, $r2) {
    for ($c = 6; $c -le 22; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-RowOdds 31 32
Swap-RowOdds 36 37
Swap-RowOdds 72 73
Swap-RowOdds 78 79
Swap-RowOdds 87 88
Swap-RowOdds 91 92
Swap-RowOdds 95 96

# Append two new match rows (106, 107) at the bottom of the table.
# Copy row 105's formatting (bold/bordered index cell style, date-time
# number format on column E) down into the new rows before setting values.
$ws.Range("A105:V105").Copy($ws.Range("A106:V106"))
$ws.Range("A105:V105").Copy($ws.Range("A107:V107"))

# Row 106
$ws.Cells.Item(106, 1).Value2 = 105
$ws.Cells.Item(106, 2).Value2 = 'netherlands'
$ws.Cells.Item(106, 3).Value2 = 'eredivisie'
$ws.Cells.Item(106, 4).Value2 = '2023-2024'
$ws.Cells.Item(106, 5).Value2 = 45242.69791666666
$ws.Cells.Item(106, 6).Value2 = 'Feyenoord'
$ws.Cells.Item(106, 7).Value2 = 1
$ws.Cells.Item(106, 8).Value2 = 'AZ Alkmaar'
$ws.Cells.Item(106, 9).Value2 = 0
$ws.Cells.Item(106, 10).Value2 = 1.67
$ws.Cells.Item(106, 11).Value2 = '04/11/2023 19:13'
$ws.Cells.Item(106, 12).Value2 = 1.5
$ws.Cells.Item(106, 13).Value2 = '12/11/2023 16:30'
$ws.Cells.Item(106, 14).Value2 = 4.22
$ws.Cells.Item(106, 15).Value2 = '04/11/2023 19:13'
$ws.Cells.Item(106, 16).Value2 = 4.71
$ws.Cells.Item(106, 17).Value2 = '12/11/2023 16:35'
$ws.Cells.Item(106, 18).Value2 = 4.82
$ws.Cells.Item(106, 19).Value2 = '04/11/2023 19:13'
$ws.Cells.Item(106, 20).Value2 = 6.66
$ws.Cells.Item(106, 21).Value2 = '12/11/2023 16:42'
$ws.Cells.Item(106, 22).Value2 = 'https://www.betexplorer.com/football/netherlands/eredivisie/feyenoord-az-alkmaar/C4oDPN4f/'

# Row 107
$ws.Cells.Item(107, 1).Value2 = 106
$ws.Cells.Item(107, 2).Value2 = 'netherlands'
$ws.Cells.Item(107, 3).Value2 = 'eredivisie'
$ws.Cells.Item(107, 4).Value2 = '2023-2024'
$ws.Cells.Item(107, 5).Value2 = 45242.83333333334
$ws.Cells.Item(107, 6).Value2 = 'Utrecht'
$ws.Cells.Item(107, 7).Value2 = 2
$ws.Cells.Item(107, 8).Value2 = 'Excelsior'
$ws.Cells.Item(107, 9).Value2 = 2
$ws.Cells.Item(107, 10).Value2 = 1.71
$ws.Cells.Item(107, 11).Value2 = '05/11/2023 12:42'
$ws.Cells.Item(107, 12).Value2 = 1.53
$ws.Cells.Item(107, 13).Value2 = '12/11/2023 19:58'
$ws.Cells.Item(107, 14).Value2 = 4.26
$ws.Cells.Item(107, 15).Value2 = '05/11/2023 12:42'
$ws.Cells.Item(107, 16).Value2 = 4.63
$ws.Cells.Item(107, 17).Value2 = '12/11/2023 19:59'
$ws.Cells.Item(107, 18).Value2 = 4.5
$ws.Cells.Item(107, 19).Value2 = '05/11/2023 12:42'
$ws.Cells.Item(107, 20).Value2 = 6.18
$ws.Cells.Item(107, 21).Value2 = '12/11/2023 19:59'
$ws.Cells.Item(107, 22).Value2 = 'https://www.betexplorer.com/football/netherlands/eredivisie/utrecht-excelsior/jepHOsK0/'

